$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values are stored as text, matching the
# original inline-string cell contents (not auto-converted to numbers).
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "295.67"
$ws.Range("E2").Value = "1.29%"
$ws.Range("G2").Value = "13"
$ws.Range("D3").Value = "40.74"
$ws.Range("E3").Value = "0.76%"
$ws.Range("G3").Value = "13"
$ws.Range("D4").Value = "5.008"
$ws.Range("E4").Value = "-0.75%"
$ws.Range("G4").Value = "13"
$ws.Range("D5").Value = "0.07385"
$ws.Range("E5").Value = "0.24%"
$ws.Range("G5").Value = "13"
$ws.Range("D6").Value = "1.555"
$ws.Range("E6").Value = "-0.25%"
$ws.Range("G6").Value = "13"
$ws.Range("D7").Value = "0.9254"
$ws.Range("E7").Value = "0.76%"
$ws.Range("G7").Value = "13"
$ws.Range("D8").Value = "2.333"
$ws.Range("E8").Value = "-2.75%"
$ws.Range("G8").Value = "13"
$ws.Range("D9").Value = "0.1194"
$ws.Range("E9").Value = "0.64%"
$ws.Range("G9").Value = "13"
$ws.Range("D10").Value = "0.1809"
$ws.Range("E10").Value = "2.63%"
$ws.Range("G10").Value = "13"
$ws.Range("D11").Value = "0.04395"
$ws.Range("E11").Value = "4.43%"
$ws.Range("G11").Value = "13"
$ws.Range("D12").Value = "0.08820"
$ws.Range("E12").Value = "0.76%"
$ws.Range("G12").Value = "13"
$ws.Range("E13").Value = "0.31%"
$ws.Range("G13").Value = "13"
$ws.Range("D14").Value = "0.001264"
$ws.Range("E14").Value = "-0.81%"
$ws.Range("G14").Value = "13"
$ws.Range("D15").Value = "0.005825"
$ws.Range("E15").Value = "-0.17%"
$ws.Range("G15").Value = "13"
$ws.Range("D16").Value = "3.359"
$ws.Range("E16").Value = "-1.60%"
$ws.Range("G16").Value = "13"
$ws.Range("E17").Value = "0.03%"
$ws.Range("G17").Value = "13"
$ws.Range("D18").Value = "0.3275"
$ws.Range("E18").Value = "-0.72%"
$ws.Range("G18").Value = "13"
$ws.Range("D19").Value = "7.889"
$ws.Range("E19").Value = "4.08%"
$ws.Range("G19").Value = "13"
$ws.Range("D20").Value = "0.1380"
$ws.Range("E20").Value = "2.76%"
$ws.Range("G20").Value = "13"
$ws.Range("D21").Value = "0.2804"
$ws.Range("E21").Value = "-2.04%"
$ws.Range("G21").Value = "13"
$ws.Range("D22").Value = "0.03923"
$ws.Range("E22").Value = "2.01%"
$ws.Range("G22").Value = "13"
$ws.Range("D23").Value = "0.001270"
$ws.Range("E23").Value = "-1.00%"
$ws.Range("G23").Value = "13"
$ws.Range("D24").Value = "0.003805"
$ws.Range("E24").Value = "-2.26%"
$ws.Range("G24").Value = "13"
$ws.Range("D25").Value = "0.0001231"
$ws.Range("E25").Value = "-4.04%"
$ws.Range("G25").Value = "13"
$ws.Range("D26").Value = "0.0003723"
$ws.Range("E26").Value = "-0.20%"
$ws.Range("G26").Value = "13"
$ws.Range("G27").Value = "13"
$ws.Range("G28").Value = "13"
$ws.Range("G29").Value = "13"
$ws.Range("G30").Value = "13"
$ws.Range("G31").Value = "13"
$ws.Range("G32").Value = "13"
$ws.Range("G33").Value = "13"
$ws.Range("G34").Value = "13"
$ws.Range("G35").Value = "13"
$ws.Range("G36").Value = "13"
$ws.Range("G37").Value = "13"
$ws.Range("D38").Value = "0.02339"
$ws.Range("E38").Value = "0.16%"
$ws.Range("G38").Value = "13"
$ws.Range("D39").Value = "0.05084"
$ws.Range("E39").Value = "1.06%"
$ws.Range("G39").Value = "13"
$ws.Range("D40").Value = "0.006176"
$ws.Range("E40").Value = "23.33%"
$ws.Range("G40").Value = "13"
$ws.Range("D41").Value = "0.007830"
$ws.Range("E41").Value = "1.84%"
$ws.Range("G41").Value = "13"
$ws.Range("D42").Value = "0.1295"
$ws.Range("E42").Value = "1.83%"
$ws.Range("G42").Value = "13"
$ws.Range("D43").Value = "0.007388"
$ws.Range("E43").Value = "0.22%"
$ws.Range("G43").Value = "13"
$ws.Range("D44").Value = "0.007364"
$ws.Range("E44").Value = "-4.14%"
$ws.Range("G44").Value = "13"
$ws.Range("D45").Value = "0.2942"
$ws.Range("E45").Value = "-7.07%"
$ws.Range("G45").Value = "13"
$ws.Range("D46").Value = "0.00006112"
$ws.Range("E46").Value = "-6.68%"
$ws.Range("G46").Value = "13"
$ws.Range("E47").Value = "-0.20%"
$ws.Range("G47").Value = "13"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "0.04661"
$ws.Range("E48").Value = "-81.49%"
$ws.Range("G48").Value = "13"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "0.004202"
$ws.Range("E49").Value = "-0.19%"
$ws.Range("G49").Value = "13"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").Value = "-0.20%"
$ws.Range("G50").Value = "13"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").Value = "-0.20%"
$ws.Range("G51").Value = "13"
